# Update the "Förändrad" (Changed) date column (C) for rows 2-43 from
# 45752 (2025-04-05) to 45753 (2025-04-06) on the "Avverkningsanmälningar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45752) {
        $cell.Value2 = 45753
    }
}
